$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 311.78946
$ws.Range("I33").Value = 264.93332
$ws.Range("J33").Value = 487.5
$ws.Range("K33").Value = 264.93332
$ws.Range("L33").Value = 487.5
$ws.Range("M33").Value = -35.93331999999998
$ws.Range("N33").Value = -945.5
$ws.Range("H95").Value = 25666.666
$ws.Range("J95").Value = 25666.666
$ws.Range("L95").Value = 25666.666
$ws.Range("N95").Value = -31158.666
$ws.Range("H98").Value = 1587.125
$ws.Range("I98").Value = 1587.125
$ws.Range("K98").Value = 1587.125
$ws.Range("M98").Value = -89.125
$ws.Range("H113").Value = 3422.1538
$ws.Range("I113").Value = 2680
$ws.Range("J113").Value = 3644.8
$ws.Range("K113").Value = 2680
$ws.Range("L113").Value = 3644.8
$ws.Range("M113").Value = 574
$ws.Range("N113").Value = -10152.8
$ws.Range("H116").Value = 2358.077
$ws.Range("I116").Value = 2031.875
$ws.Range("J116").Value = 2880
$ws.Range("K116").Value = 2031.875
$ws.Range("L116").Value = 2880
$ws.Range("M116").Value = 1410.125
$ws.Range("N116").Value = -9764
$ws.Range("H122").Value = 1587.125
$ws.Range("I122").Value = 1587.125
$ws.Range("K122").Value = 4761.375
$ws.Range("M122").Value = -2311.375
$ws.Range("H137").Value = 536308.4
$ws.Range("I137").Value = 2036.4166
$ws.Range("K137").Value = 6109.2498
$ws.Range("M137").Value = -3559.2498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21003
$ws.Range("I32").Value = 22711.02
$ws.Range("J32").Value = 12121.3
$ws.Range("K32").Value = 22711.02
$ws.Range("L32").Value = 12121.3
$ws.Range("M32").Value = -22424.02
$ws.Range("N32").Value = -12695.3
$ws.Range("H61").Value = 5897.881
$ws.Range("I61").Value = 3903.1936
$ws.Range("J61").Value = 11519.272
$ws.Range("K61").Value = 3903.1936
$ws.Range("L61").Value = 11519.272
$ws.Range("M61").Value = -3691.1936
$ws.Range("N61").Value = -11943.272
$ws.Range("H74").Value = 3632.149
$ws.Range("I74").Value = 1601.2812
$ws.Range("J74").Value = 7964.6665
$ws.Range("K74").Value = 1601.2812
$ws.Range("L74").Value = 7964.6665
$ws.Range("M74").Value = -727.2811999999999
$ws.Range("N74").Value = -9712.666499999999
$ws.Range("H77").Value = 3632.149
$ws.Range("I77").Value = 1601.2812
$ws.Range("J77").Value = 7964.6665
$ws.Range("K77").Value = 8006.405999999999
$ws.Range("L77").Value = 39823.3325
$ws.Range("M77").Value = -3638.405999999999
$ws.Range("N77").Value = -48559.3325
$ws.Range("H80").Value = 27277.25
$ws.Range("J80").Value = 40054.5
$ws.Range("L80").Value = 40054.5
$ws.Range("N80").Value = -42050.5
$ws.Range("H83").Value = 27277.25
$ws.Range("J83").Value = 40054.5
$ws.Range("L83").Value = 120163.5
$ws.Range("N83").Value = -130147.5
$ws.Range("H102").Value = 3360
$ws.Range("I102").Value = 2977.7778
$ws.Range("J102").Value = 3672.7273
$ws.Range("K102").Value = 2977.7778
$ws.Range("L102").Value = 3672.7273
$ws.Range("M102").Value = -1355.7778
$ws.Range("N102").Value = -6916.7273
$ws.Range("H109").Value = 48900
$ws.Range("J109").Value = 48900
$ws.Range("L109").Value = 48900
$ws.Range("N109").Value = -51674
$ws.Range("H132").Value = 4072.0942
$ws.Range("I132").Value = 1454.5938
$ws.Range("J132").Value = 8060.6665
$ws.Range("K132").Value = 4363.7814
$ws.Range("L132").Value = 24181.9995
$ws.Range("M132").Value = -1833.7814
$ws.Range("N132").Value = -29241.9995
$ws.Range("H136").Value = 5897.881
$ws.Range("I136").Value = 3903.1936
$ws.Range("J136").Value = 11519.272
$ws.Range("K136").Value = 11709.5808
$ws.Range("L136").Value = 34557.81600000001
$ws.Range("M136").Value = -9159.5808
$ws.Range("N136").Value = -39657.81600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 500
$ws.Range("I64").Value = 500
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 500
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -275
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 500
$ws.Range("I67").Value = 500
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 500
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 280
$ws.Range("N67").ClearContents()
$ws.Range("H99").Value = 749.5
$ws.Range("I99").Value = 832.6667
$ws.Range("J99").Value = 500
$ws.Range("K99").Value = 832.6667
$ws.Range("L99").Value = 500
$ws.Range("M99").Value = 665.3333
$ws.Range("N99").Value = -3496
$ws.Range("H134").Value = 2357.258
$ws.Range("I134").Value = 2383.476
$ws.Range("J134").Value = 2302.2
$ws.Range("K134").Value = 7150.428
$ws.Range("L134").Value = 6906.599999999999
$ws.Range("M134").Value = -4615.428
$ws.Range("N134").Value = -11976.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5553
$ws.Range("I31").Value = 8007.6875
$ws.Range("J31").Value = 3485.8948
$ws.Range("K31").Value = 8007.6875
$ws.Range("L31").Value = 3485.8948
$ws.Range("M31").Value = -7712.6875
$ws.Range("N31").Value = -4075.8948
$ws.Range("H34").Value = 5553
$ws.Range("I34").Value = 8007.6875
$ws.Range("J34").Value = 3485.8948
$ws.Range("K34").Value = 8007.6875
$ws.Range("L34").Value = 3485.8948
$ws.Range("M34").Value = -7805.6875
$ws.Range("N34").Value = -3889.8948
$ws.Range("H105").Value = 2132.6667
$ws.Range("I105").Value = 899
$ws.Range("J105").Value = 2749.5
$ws.Range("K105").Value = 899
$ws.Range("L105").Value = 2749.5
$ws.Range("M105").Value = 848
$ws.Range("N105").Value = -6243.5
$ws.Range("H132").Value = 1885.1666
$ws.Range("I132").Value = 1471.2894
$ws.Range("J132").Value = 2868.125
$ws.Range("K132").Value = 4413.8682
$ws.Range("L132").Value = 8604.375
$ws.Range("M132").Value = -1883.8682
$ws.Range("N132").Value = -13664.375
$ws.Range("H134").Value = 2499.1453
$ws.Range("I134").Value = 1383.0625
$ws.Range("J134").Value = 3689.6333
$ws.Range("K134").Value = 4149.1875
$ws.Range("L134").Value = 11068.8999
$ws.Range("M134").Value = -1614.1875
$ws.Range("N134").Value = -16138.8999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 42542.5
$ws.Range("J62").Value = 42542.5
$ws.Range("L62").Value = 42542.5
$ws.Range("N62").Value = -43914.5
$ws.Range("H65").Value = 42542.5
$ws.Range("J65").Value = 42542.5
$ws.Range("L65").Value = 127627.5
$ws.Range("N65").Value = -134491.5
$ws.Range("H70").Value = 5793.619
$ws.Range("I70").Value = 5130.857
$ws.Range("K70").Value = 5130.857
$ws.Range("M70").Value = -4860.857
$ws.Range("H73").Value = 5793.619
$ws.Range("I73").Value = 5130.857
$ws.Range("K73").Value = 5130.857
$ws.Range("M73").Value = -4194.857
$ws.Range("H122").Value = 34489.668
$ws.Range("I122").Value = 49999.5
$ws.Range("J122").Value = 3470
$ws.Range("K122").Value = 149998.5
$ws.Range("L122").Value = 10410
$ws.Range("M122").Value = -147548.5
$ws.Range("N122").Value = -15310
$ws.Range("H132").Value = 2569.5862
$ws.Range("I132").Value = 2375.3333
$ws.Range("K132").Value = 7125.999899999999
$ws.Range("M132").Value = -4595.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3461.5386
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 3750
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 3750
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -4022
$ws.Range("H68").Value = 2500
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251
$ws.Range("H71").Value = 2500
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256
$ws.Range("H74").Value = 33478
$ws.Range("I74").Value = 28000
$ws.Range("K74").Value = 28000
$ws.Range("M74").Value = -27002
$ws.Range("H77").Value = 33478
$ws.Range("I77").Value = 28000
$ws.Range("K77").Value = 84000
$ws.Range("M77").Value = -79008
$ws.Range("H93").Value = 302.6
$ws.Range("I93").Value = 211.58333
$ws.Range("J93").Value = 666.6667
$ws.Range("K93").Value = 211.58333
$ws.Range("L93").Value = 666.6667
$ws.Range("M93").Value = 1036.41667
$ws.Range("N93").Value = -3162.6667
$ws.Range("H132").Value = 6887.7744
$ws.Range("I132").Value = 10337.866
$ws.Range("J132").Value = 3653.3125
$ws.Range("K132").Value = 31013.598
$ws.Range("L132").Value = 10959.9375
$ws.Range("M132").Value = -28483.598
$ws.Range("N132").Value = -16019.9375
$ws.Range("H136").Value = 5209.561
$ws.Range("I136").Value = 2996.4546
$ws.Range("J136").Value = 7772.1055
$ws.Range("K136").Value = 8989.363799999999
$ws.Range("L136").Value = 23316.3165
$ws.Range("M136").Value = -6439.363799999999
$ws.Range("N136").Value = -28416.3165

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H107").Value = 1982.2727
$ws.Range("I107").Value = 665.6667
$ws.Range("J107").Value = 2476
$ws.Range("K107").Value = 1997.0001
$ws.Range("L107").Value = 7428
$ws.Range("M107").Value = -77.00009999999997
$ws.Range("N107").Value = -11268
$ws.Range("H123").Value = 24296
$ws.Range("J123").Value = 24296
$ws.Range("L123").Value = 24296
$ws.Range("N123").Value = -34096
$ws.Range("H132").Value = 1440.566
$ws.Range("I132").Value = 694.2143
$ws.Range("K132").Value = 2082.6429
$ws.Range("M132").Value = 447.3571000000002
$ws.Range("H136").Value = 6028.593
$ws.Range("I136").Value = 4152.324
$ws.Range("J136").Value = 10112.235
$ws.Range("K136").Value = 12456.972
$ws.Range("L136").Value = 30336.705
$ws.Range("M136").Value = -9906.971999999998
$ws.Range("N136").Value = -35436.705

